$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.306.47'
$ws.Range('E2').Value = '  +0.63%  '
$ws.Range('D3').Value = '3.740.69'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '616.48'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +5.49%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '187.35'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +5.30%  '
$ws.Range('D7').Value = '3.738.64'
$ws.Range('E7').Value = '  +0.43%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.642'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.43%  '
$ws.Range('E9').Value = '  -0.47%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.723'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.163'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -3.58%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '57.24'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +5.41%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000296'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -4.15%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '10.74'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.19%  '
$ws.Range('D15').Value = '4.336.23'
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('D16').Value = '3.744.10'
$ws.Range('E16').Value = '  -0.68%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '13.13'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '19.41'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.55%  '
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('E20').Value = '  -1.57%  '
$ws.Range('D21').Value = '69.154.45'
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '415.15'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '89.62'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.07%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.07'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.95%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '12.91'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.63%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.02'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.33%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.08'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.24%  '
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '9.72'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.79%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '33.34'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.47%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.40'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -11.95%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '12.79'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.66%  '
$ws.Range('E34').Value = '  +1.97%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '44.95'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.89%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '66.33'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '617.80'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.63%  '
$ws.Range('D38').Value = '0.0₃0871'
$ws.Range('E38').Value = '  -8.22%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.410'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.10%  '
$ws.Range('E40').Value = '  +0.20%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = "Normal"
$ws.Range('E42').Value = '  +0.88%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.08'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.58%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0446'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.80%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.66'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.142'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.36%  '
$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.28'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.32%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.853.82'
$ws.Range('E48').Value = '  +2.45%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.73'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +2.66%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.70'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -17.18%  '
$ws.Range('E51').Value = '  -3.44%  '
